# Rename Sheet1 -> Transcriptomes, add new Genomes sheet, and populate new rows
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Transcriptomes"

# New sample rows for the Transcriptomes sheet (columns A, C, E)
$newRows = @(
    ,@(43, "Trematomus bernacchii", "liver", "SRR7164558")
    ,@(44, "Trematomus bernacchii", "heart", "SRR7164559")
    ,@(45, "Trematomus bernacchii", "head kidney", "SRR7164560")
    ,@(46, "Trematomus bernacchii", "stomach", "SRR7164564")
    ,@(47, "Trematomus bernacchii", "muscle", "SRR7164565")
    ,@(48, "Trematomus bernacchii", "spleen", "SRR7164570")
    ,@(49, "Chaenocephalus aceratus", "heart", "SRR6929341")
    ,@(50, "Chaenocephalus aceratus", "gill", "SRR6929342")
    ,@(51, "Chaenocephalus aceratus", "muscle", "SRR6929345")
    ,@(52, "Chaenocephalus aceratus", "liver", "SRR6929346")
    ,@(53, "Chaenocephalus aceratus", "kidney", "SRR6929347")
    ,@(54, "Chaenocephalus aceratus", "intestine", "SRR6929348")
    ,@(55, "Eleginops maclovinus", "liver", "SRR6793933")
    ,@(56, "Eleginops maclovinus", "red muscle", "SRR6793934")
    ,@(57, "Eleginops maclovinus", "gill", "SRR6793935")
    ,@(58, "Eleginops maclovinus", "small intestine", "SRR6793936")
    ,@(59, "Eleginops maclovinus", "head kidney", "SRR6793937")
    ,@(60, "Eleginops maclovinus", "spleen", "SRR6793939")
    ,@(61, "Gymnodraco acuticeps", "adult tissue pool (brain, gill, liver, spleen)", "SRR6450838")
    ,@(62, "Pagothenia borchgrevinki", "ref transcriptome from multiple tissues", "SRR5210464")
    ,@(63, "Pagothenia borchgrevinki", "heat-stressed gill samples", "SRR5210375")
    ,@(64, "Chionodraco rastrospinosus", "ref transcriptome from multiple tissues", "SRR5210463")
    ,@(65, "Chionodraco rastrospinosus", "heat-stressed gill samples", "SRR5210373")
    ,@(66, "Chionodraco hamatus", "gill transcriptome", "SRR4279902")
    ,@(67, "Notothenia coriiceps", "adult Poly I:C challenged liver", "SRR3342843")
    ,@(68, "Notothenia coriiceps", "adult HKEB challenged liver", "SRR3342842")
    ,@(69, "Notothenia coriiceps", "adult control liver", "SRR3342841")
    ,@(70, "Notothenia coriiceps", "pronephric kidney", "SRR3133082")
    ,@(71, "Notothenia coriiceps", "pronephric kidney", "SRR3133083")
    ,@(72, "Trematomus pennellii", "brain", "SRR2822458")
    ,@(73, "Trematomus pennellii", "liver", "SRR2823736")
    ,@(74, "Trematomus newnesi", "pooled PE reads", "SRR2259813")
    ,@(75, "Trematomus newnesi", "pooled PE reads", "SRR2259814")
    ,@(76, "Lepidonotothen nudifrons", "pooled spleen samples from different temperature treatments", "ERR793598")
    ,@(77, "Lepidonotothen nudifrons", "pooled spleen samples from different temperature treatments", "ERR793597")
    ,@(78, "Gymnodraco acuticeps", "?", "SRR2072641")
    ,@(79, "Gymnodraco acuticeps", "?", "SRR2072640")
    ,@(80, "Gymnodraco acuticeps", "?", "SRR2072639")
    ,@(81, "Chionodraco hamatus", "?", "SRR2072638")
    ,@(82, "Chionodraco hamatus", "?", "SRR2072637")
    ,@(83, "Chionodraco hamatus", "?", "SRR2072636")
    ,@(84, "Dissostichus mawsoni", "head kidney", "SRR6794059")
    ,@(85, "Dissostichus mawsoni", "brain", "SRR6794060")
    ,@(86, "Dissostichus mawsoni", "liver", "SRR6794061")
    ,@(87, "Dissostichus mawsoni", "gill", "SRR6794062")
    ,@(88, "Dissostichus mawsoni", "red muscle", "SRR6794063")
    ,@(89, "Dissostichus mawsoni", "white muscle", "SRR6794064")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $eCell = $ws.Cells.Item($r, 5)
    $eCell.Value = $row[3]
    $eCell.Font.Color = 255
}

# Add the new Genomes worksheet after Transcriptomes
$ws2 = $wb.Worksheets.Add($null, $ws)
$ws2.Name = "Genomes"

$genomeRows = @(
    ,@(1, "Chaenocephalus aceratus", "SRR6942631")
    ,@(2, "Chaenocephalus aceratus", "SRR6942632")
    ,@(4, "Eleginops maclovinus", $null)
)

foreach ($row in $genomeRows) {
    $r = $row[0]
    $ws2.Cells.Item($r, 1).Value = $row[1]
    if ($row[2] -ne $null) {
        $ws2.Cells.Item($r, 2).Value = $row[2]
    }
}

$ws2.Columns.Item(1).AutoFit()
$ws2.Range("B4").Select()

